$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "How many curves can I load in one go?"
$ws.Range("B5").Value = "llama3.2:latest"
$ws.Range("C5").Value = "You can load up to 450 curves at a time."
